# Updates the cryptocurrency price table (columns B-E) to the latest
# scraped values, as produced by the scheduled GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.712.77"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "'2.200.65"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'229.57"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").Value = "'60.37"
$ws.Range("E7").Value = "  -2.26%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "'0.400"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "'56.87"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").Value = "'0.0890"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "'2.528.97"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "'15.33"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "'22.01"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "'0.793"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "'5.55"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "'2.200.57"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "'41.651.09"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "'0.0₃0906"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").Value = "'71.80"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "'240.68"
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'2.34"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").Value = "'2.27"
$ws.Range("E26").Value = "  -5.17%  "
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").Value = "'167.78"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "'0.139"
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "'19.68"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").Value = "'2.60"
$ws.Range("E32").Value = "  -7.14%  "
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("D34").Value = "'4.94"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").Value = "'4.58"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").Value = "'0.0643"
$ws.Range("E36").Value = "  +3.17%  "
$ws.Range("D37").Value = "'6.30"
$ws.Range("E37").Value = "  -5.81%  "
$ws.Range("E38").Value = "  -6.22%  "
$ws.Range("D39").Value = "'2.33"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'0.0239"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").Value = "'8.57"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'1.20"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0951"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("D46").Value = "'96.34"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("E47").Value = "  -11.20%  "
$ws.Range("D48").Value = "'1.448.76"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").Value = "'16.02"
$ws.Range("E50").Value = "  -3.67%  "
$ws.Range("E51").Value = "  -1.19%  "
